$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'257.82"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = "'0.62%"
$ws.Range('E2').Style = 'Normal'
$ws.Range('G2').Value = "'21"
$ws.Range('G2').Style = 'Normal'
$ws.Range('D3').Value = "'27.19"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = "'-1.39%"
$ws.Range('E3').Style = 'Normal'
$ws.Range('G3').Value = "'21"
$ws.Range('G3').Style = 'Normal'
$ws.Range('D4').Value = "'4.663"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = "'-10.63%"
$ws.Range('E4').Style = 'Normal'
$ws.Range('G4').Value = "'21"
$ws.Range('G4').Style = 'Normal'
$ws.Range('D5').Value = "'0.05887"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = "'-0.53%"
$ws.Range('E5').Style = 'Normal'
$ws.Range('G5').Value = "'21"
$ws.Range('G5').Style = 'Normal'
$ws.Range('D6').Value = "'6.645"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = "'-0.41%"
$ws.Range('E6').Style = 'Normal'
$ws.Range('G6').Value = "'21"
$ws.Range('G6').Style = 'Normal'
$ws.Range('E7').Value = "'-0.97%"
$ws.Range('E7').Style = 'Normal'
$ws.Range('G7').Value = "'21"
$ws.Range('G7').Style = 'Normal'
$ws.Range('D8').Value = "'0.9428"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = "'-8.64%"
$ws.Range('E8').Style = 'Normal'
$ws.Range('G8').Value = "'21"
$ws.Range('G8').Style = 'Normal'
$ws.Range('D9').Value = "'0.1406"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = "'-0.90%"
$ws.Range('E9').Style = 'Normal'
$ws.Range('G9').Value = "'21"
$ws.Range('G9').Style = 'Normal'
$ws.Range('D10').Value = "'0.03860"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = "'6.26%"
$ws.Range('E10').Style = 'Normal'
$ws.Range('G10').Value = "'21"
$ws.Range('G10').Style = 'Normal'
$ws.Range('D11').Value = "'0.07098"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = "'-1.30%"
$ws.Range('E11').Style = 'Normal'
$ws.Range('G11').Value = "'21"
$ws.Range('G11').Style = 'Normal'
$ws.Range('D12').Value = "'0.03179"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = "'-1.95%"
$ws.Range('E12').Style = 'Normal'
$ws.Range('G12').Value = "'21"
$ws.Range('G12').Style = 'Normal'
$ws.Range('D13').Value = "'0.09184"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = "'-0.46%"
$ws.Range('E13').Style = 'Normal'
$ws.Range('G13').Value = "'21"
$ws.Range('G13').Style = 'Normal'
$ws.Range('D14').Value = "'0.001543"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = "'0.17%"
$ws.Range('E14').Style = 'Normal'
$ws.Range('G14').Value = "'21"
$ws.Range('G14').Style = 'Normal'
$ws.Range('B15').Value = "'TigerCash"
$ws.Range('B15').Style = 'Normal'
$ws.Range('C15').Value = "'https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range('C15').Style = 'Normal'
$ws.Range('D15').Value = "'0.006220"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = "'8.18%"
$ws.Range('E15').Style = 'Normal'
$ws.Range('G15').Value = "'21"
$ws.Range('G15').Style = 'Normal'
$ws.Range('B16').Value = "'LEO"
$ws.Range('B16').Style = 'Normal'
$ws.Range('C16').Value = "'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range('C16').Style = 'Normal'
$ws.Range('D16').Value = "'3.514"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = "'1.10%"
$ws.Range('E16').Style = 'Normal'
$ws.Range('G16').Value = "'21"
$ws.Range('G16').Style = 'Normal'
$ws.Range('B17').Value = "'GateToken"
$ws.Range('B17').Style = 'Normal'
$ws.Range('C17').Value = "'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range('C17').Style = 'Normal'
$ws.Range('D17').Value = "'3.207"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = "'-1.75%"
$ws.Range('E17').Style = 'Normal'
$ws.Range('G17').Value = "'21"
$ws.Range('G17').Style = 'Normal'
$ws.Range('B18').Value = "'BTSEToken"
$ws.Range('B18').Style = 'Normal'
$ws.Range('C18').Value = "'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range('C18').Style = 'Normal'
$ws.Range('D18').Value = "'2.225"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = "'0.09%"
$ws.Range('E18').Style = 'Normal'
$ws.Range('G18').Value = "'21"
$ws.Range('G18').Style = 'Normal'
$ws.Range('B19').Value = "'One"
$ws.Range('B19').Style = 'Normal'
$ws.Range('C19').Value = "'https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range('C19').Style = 'Normal'
$ws.Range('D19').Value = "'0.0006039"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = "'-0.66%"
$ws.Range('E19').Style = 'Normal'
$ws.Range('G19').Value = "'21"
$ws.Range('G19').Style = 'Normal'
$ws.Range('G20').Value = "'21"
$ws.Range('G20').Style = 'Normal'
$ws.Range('D21').Value = "'0.1292"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = "'-1.26%"
$ws.Range('E21').Style = 'Normal'
$ws.Range('G21').Value = "'21"
$ws.Range('G21').Style = 'Normal'
$ws.Range('D22').Value = "'3.877"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = "'9.82%"
$ws.Range('E22').Style = 'Normal'
$ws.Range('G22').Value = "'21"
$ws.Range('G22').Style = 'Normal'
$ws.Range('D23').Value = "'0.04225"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = "'1.23%"
$ws.Range('E23').Style = 'Normal'
$ws.Range('G23').Value = "'21"
$ws.Range('G23').Style = 'Normal'
$ws.Range('D24').Value = "'0.001219"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = "'-0.06%"
$ws.Range('E24').Style = 'Normal'
$ws.Range('G24').Value = "'21"
$ws.Range('G24').Style = 'Normal'
$ws.Range('E25').Value = "'-4.77%"
$ws.Range('E25').Style = 'Normal'
$ws.Range('G25').Value = "'21"
$ws.Range('G25').Style = 'Normal'
$ws.Range('D26').Value = "'0.0001200"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = "'-0.03%"
$ws.Range('E26').Style = 'Normal'
$ws.Range('G26').Value = "'21"
$ws.Range('G26').Style = 'Normal'
$ws.Range('D27').Value = "'0.0001937"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = "'-0.23%"
$ws.Range('E27').Style = 'Normal'
$ws.Range('G27').Value = "'21"
$ws.Range('G27').Style = 'Normal'
$ws.Range('G28').Value = "'21"
$ws.Range('G28').Style = 'Normal'
$ws.Range('G29').Value = "'21"
$ws.Range('G29').Style = 'Normal'
$ws.Range('G30').Value = "'21"
$ws.Range('G30').Style = 'Normal'
$ws.Range('G31').Value = "'21"
$ws.Range('G31').Style = 'Normal'
$ws.Range('G32').Value = "'21"
$ws.Range('G32').Style = 'Normal'
$ws.Range('G33').Value = "'21"
$ws.Range('G33').Style = 'Normal'
$ws.Range('G34').Value = "'21"
$ws.Range('G34').Style = 'Normal'
$ws.Range('G35').Value = "'21"
$ws.Range('G35').Style = 'Normal'
$ws.Range('G36').Value = "'21"
$ws.Range('G36').Style = 'Normal'
$ws.Range('G37').Value = "'21"
$ws.Range('G37').Style = 'Normal'
$ws.Range('G38').Value = "'21"
$ws.Range('G38').Style = 'Normal'
$ws.Range('G39').Value = "'21"
$ws.Range('G39').Style = 'Normal'
$ws.Range('D40').Value = "'0.03832"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = "'0.49%"
$ws.Range('E40').Style = 'Normal'
$ws.Range('G40').Value = "'21"
$ws.Range('G40').Style = 'Normal'
$ws.Range('D41').Value = "'0.006230"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = "'57.07%"
$ws.Range('E41').Style = 'Normal'
$ws.Range('G41').Value = "'21"
$ws.Range('G41').Style = 'Normal'
$ws.Range('D42').Value = "'0.1102"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = "'-0.19%"
$ws.Range('E42').Style = 'Normal'
$ws.Range('G42').Value = "'21"
$ws.Range('G42').Style = 'Normal'
$ws.Range('D43').Value = "'0.002200"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = "'-10.60%"
$ws.Range('E43').Style = 'Normal'
$ws.Range('G43').Value = "'21"
$ws.Range('G43').Style = 'Normal'
$ws.Range('E44').Value = "'7.19%"
$ws.Range('E44').Style = 'Normal'
$ws.Range('G44').Value = "'21"
$ws.Range('G44').Style = 'Normal'
$ws.Range('D45').Value = "'0.00005457"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = "'0.38%"
$ws.Range('E45').Style = 'Normal'
$ws.Range('G45').Value = "'21"
$ws.Range('G45').Style = 'Normal'
$ws.Range('E46').Value = "'-0.18%"
$ws.Range('E46').Style = 'Normal'
$ws.Range('G46').Value = "'21"
$ws.Range('G46').Style = 'Normal'
$ws.Range('E47').Value = "'-45.10%"
$ws.Range('E47').Style = 'Normal'
$ws.Range('G47').Value = "'21"
$ws.Range('G47').Style = 'Normal'
$ws.Range('D48').Value = "'0.1314"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = "'5,956.94%"
$ws.Range('E48').Style = 'Normal'
$ws.Range('G48').Value = "'21"
$ws.Range('G48').Style = 'Normal'
$ws.Range('D49').Value = "'0.00002100"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = "'-0.18%"
$ws.Range('E49').Style = 'Normal'
$ws.Range('G49').Value = "'21"
$ws.Range('G49').Style = 'Normal'
$ws.Range('D50').Value = "'0.0002000"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = "'-0.18%"
$ws.Range('E50').Style = 'Normal'
$ws.Range('G50').Value = "'21"
$ws.Range('G50').Style = 'Normal'
$ws.Range('G51').Value = "'21"
$ws.Range('G51').Style = 'Normal'

Write-Host "Edit complete"
